# Agregando nuevas búsquedas ajax
# Insert a new transaction row at the top of the statement (row 1), pushing
# all existing rows down by one. The new row represents the most recent
# account movement (2014-03-25, document 0000950793, balance 4112.84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every existing row down by one, opening up a blank row 1.
$ws.Rows("1:1").Insert()

# Copy the (now shifted) second row's formatting onto the new first row so
# the date / text / amount columns keep their original number formats
# without Excel creating brand-new style entries.
$ws.Range("A2:H2").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)

# Non-breaking spaces used by the bank export after the amount value.
$nbsp = [char]160
$monto = "0.26" + "$nbsp$nbsp"

$ws.Range("A1").Value = 41723
$ws.Range("B1").Value = "INTERES A SU FAVOR"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "0000950793"
$ws.Range("E1").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F1").Value = $monto
$ws.Range("G1").Value = "4112.84"

# Re-create the CONCATENATE "export" formula on the new row 1 (adapted from
# the shifted copy that now lives in row 2), then remove the stray copy the
# row-insert left behind in H2 - only row 1 ever carries this helper formula.
$f = $ws.Range("H2").Formula
$f1 = $f -replace 'A2','A1' -replace 'B2','B1' -replace 'C2','C1' -replace 'D2','D1' -replace 'E2','E1' -replace 'F2','F1' -replace 'G2','G1'
$ws.Range("H1").Formula = $f1
$ws.Range("H2").ClearContents()

Write-Output "Inserted new row 1 for document 0000950793 / balance 4112.84"
